$wb = $excel.ActiveWorkbook

# Rows requiring updates on both the "展览" and "全部类型" sheets (identical F-column values)
$updates = @{
    2  = 220
    4  = 12942
    5  = 1336
    6  = 206
    8  = 95
    9  = 166
    10 = 221
    11 = 467
    12 = 2
    13 = 67
    17 = 408
    18 = 5511
    20 = 51
    21 = 959
    22 = 31
    24 = 130
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
